$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-25 Thursday" "2024-07-26 Friday"

Replace-Text "363÷2=181, 1" "994÷3=331, 1"
Replace-Text "369÷7=52, 5" "226÷8=28, 2"
Replace-Text "134÷2=67, 0" "459÷6=76, 3"
Replace-Text "832÷7=118, 6" "483÷3=161, 0"
Replace-Text "659÷4=164, 3" "424÷3=141, 1"
Replace-Text "348÷9=38, 6" "348÷3=116, 0"
Replace-Text "654÷6=109, 0" "422÷3=140, 2"
Replace-Text "929÷3=309, 2" "232÷6=38, 4"
Replace-Text "674÷8=84, 2" "186÷9=20, 6"
Replace-Text "230÷7=32, 6" "550÷3=183, 1"
Replace-Text "875÷2=437, 1" "564÷8=70, 4"
Replace-Text "268÷7=38, 2" "177÷6=29, 3"
Replace-Text "231÷8=28, 7" "310÷4=77, 2"
Replace-Text "345÷8=43, 1" "500÷6=83, 2"
Replace-Text "355÷6=59, 1" "610÷4=152, 2"
Replace-Text "499÷4=124, 3" "501÷7=71, 4"
Replace-Text "490÷3=163, 1" "156÷3=52, 0"
Replace-Text "120÷6=20, 0" "533÷5=106, 3"
Replace-Text "598÷8=74, 6" "448÷5=89, 3"
Replace-Text "604÷7=86, 2" "206÷6=34, 2"
Replace-Text "848÷7=121, 1" "874÷8=109, 2"
Replace-Text "760÷6=126, 4" "759÷2=379, 1"
Replace-Text "864÷4=216, 0" "968÷4=242, 0"
Replace-Text "109÷6=18, 1" "650÷7=92, 6"
Replace-Text "562÷6=93, 4" "348÷9=38, 6"
